$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Q6)
$ws.Range("B7").Value = 0.5408191021347493
$ws.Range("C7").Value = 2.531376869170928
$ws.Range("D7").Value = 24.15856235852996
$ws.Range("E7").Value = 4.915136046797683
$ws.Range("F7").Value = 4.954590293174489
$ws.Range("G7").Value = 36

# Row 8 (Q7)
$ws.Range("B8").Value = 0.5926037554034844
$ws.Range("C8").Value = 2.592646376122079
$ws.Range("D8").Value = 25.19591940884908
$ws.Range("E8").Value = 5.019553706142517
$ws.Range("F8").Value = 5.057219379343858
$ws.Range("G8").Value = 35

# Row 9 (Q8)
$ws.Range("B9").Value = 0.6352104498595522
$ws.Range("C9").Value = 3.485517187853179
$ws.Range("D9").Value = 41.44872118360328
$ws.Range("E9").Value = 6.438068125113564
$ws.Range("F9").Value = 6.573089385636792
$ws.Range("G9").Value = 20

# Row 10 (Q9)
$ws.Range("B10").Value = -1.161731967122482
$ws.Range("C10").Value = 3.357451723839908
$ws.Range("D10").Value = 28.51140454140828
$ws.Range("E10").Value = 5.339607152348221
$ws.Range("F10").Value = 5.424505998350313
$ws.Range("G10").Value = 13

# Row 11 (Q10)
$ws.Range("B11").Value = -0.3568048824601059
$ws.Range("C11").Value = 2.794051366433416
$ws.Range("D11").Value = 9.435265004421183
$ws.Range("E11").Value = 3.07168764760045
$ws.Range("F11").Value = 3.411003386152271
$ws.Range("G11").Value = 5
